# Helper to build the OLE/VBA-style BGR-packed long from an RGB hex triplet,
# since Shape.Fill.ForeColor.RGB expects the classic VBA RGB() encoding.
function RGBVal([int]$r, [int]$g, [int]$b) {
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# "Table 5" is the 5th shape on the slide (Title, Subtitle, Table 3, Table 4, Table 5).
$sh = $s.Shapes.Item(5)
$tbl = $sh.Table

# Narrow the date column (and therefore the whole table/graphic frame) - the
# table grid's second column shrinks from 2503000 EMU to 1360000 EMU, which
# also shrinks the containing graphicFrame's cx from 2974000 to 1831000 EMU.
# PowerPoint table/shape dimensions are expressed in points (914400 EMU/in,
# 72 pt/in), so 1360000 EMU == 107.0866141732283 pt.
$tbl.Columns.Item(2).Width = 107.0866141732283

$fillColor = RGBVal 0xFF 0xF0 0xC9

# New date strings (switch from "YYYY-MM-DD 00:00:00" datetime dtype to a
# plain "YYYY/MM/DD" display string) plus a light fill on every cell of the
# 3x2 table.
$dates = @("2021/03/24", "2021/04/19", "2021/04/19")

for ($row = 1; $row -le $tbl.Rows.Count; $row++) {
    for ($col = 1; $col -le $tbl.Columns.Count; $col++) {
        $cell = $tbl.Cell($row, $col)
        $cell.Shape.Fill.ForeColor.RGB = $fillColor
    }
    $dateCell = $tbl.Cell($row, 2)
    $dateCell.Shape.TextFrame.TextRange.Text = $dates[$row - 1]
}
